# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) slide master -> drives the
#                             "Integral" / Red Violet look seen on every slide
#   ppt/theme/theme2.xml  -> bound only to the notes master
#
# The target edit swaps the two themes' content: the slide master's theme
# becomes the plain default "Office Theme" colour scheme (previously living
# in theme2.xml), while the notes-master-only theme would become the old
# "Integral" / Red Violet scheme.
#
# The PowerPoint object model only exposes an editable theme/colour scheme
# off the slide master (SlideMaster.Theme.ThemeColorScheme / Master.Design /
# Presentation.Designs(1) — NotesMaster.Theme resolves to that very same
# object), so the reachable, user-visible half of this change is recolouring
# the slide master's theme to the standard Office palette. Do that by pushing
# each of the twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) to the stock "Office Theme" RGB values, in slot order.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$officeThemeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $officeThemeRgb.Length; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeRgb[$i - 1]
}
